# Update Fecha (D), Volumen (J), Precio mínimo (K), Precio máximo (L),
# Precio promedio ponderado (M), Origen (O) and Precio $/Kg (P) values
# for the "Hortaliza, Mercado Mayorista Lo Valledor de Santiago - Achicoria"
# sheet so the weekly data rows line up with their correct dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44204
$ws.Range("J2").Value = 430

$ws.Range("D3").Value = 44208
$ws.Range("J3").Value = 160

$ws.Range("D4").Value = 44188
$ws.Range("J4").Value = 210

$ws.Range("D5").Value = 44232

$ws.Range("D6").Value = 44189
$ws.Range("J6").Value = 250

$ws.Range("D8").Value = 44230
$ws.Range("J8").Value = 250
$ws.Range("L8").Value = 6000
$ws.Range("M8").Value = 5500
$ws.Range("O8").Value = "Provincia de Quillota"
$ws.Range("P8").Value = 344

$ws.Range("D9").Value = 44187
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 5000
$ws.Range("M9").Value = 5500
$ws.Range("O9").Value = "Provincia de Quillota"
$ws.Range("P9").Value = 344

$ws.Range("D10").Value = 44186
$ws.Range("J10").Value = 160

$ws.Range("D11").Value = 44292
$ws.Range("J11").Value = 90
$ws.Range("K11").Value = 6000
$ws.Range("M11").Value = 6000
$ws.Range("O11").Value = "Región Metropolitana"
$ws.Range("P11").Value = 375

$ws.Range("D12").Value = 44215

$ws.Range("D13").Value = 44210
$ws.Range("J13").Value = 340

$ws.Range("D14").Value = 44251
$ws.Range("J14").Value = 120
$ws.Range("L14").Value = 5000
$ws.Range("M14").Value = 5000
$ws.Range("O14").Value = "Región Metropolitana"
$ws.Range("P14").Value = 312
